$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.948.82'
$ws.Range('E2').Value = '  -3.19%  '
$ws.Range('D3').Value = '1.831.39'
$ws.Range('E3').Value = '  -2.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '278.08'
$ws.Range('E5').Value = '  -7.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9985'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5113'
$ws.Range('E7').Value = '  -4.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3484'
$ws.Range('E8').Value = '  -6.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.67'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06799'
$ws.Range('E10').Value = '  -4.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.85'
$ws.Range('E11').Value = '  -7.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8043'
$ws.Range('E12').Value = '  -9.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07805'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('D14').Value = '1.829.84'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.069'
$ws.Range('E15').Value = '  -4.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.14'
$ws.Range('E16').Value = '  -4.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9982'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.15'
$ws.Range('E18').Value = '  -4.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008061'
$ws.Range('E19').Value = '  -4.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9982'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = '25.983.19'
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.763'
$ws.Range('E22').Value = '  -4.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.03'
$ws.Range('E23').Value = '  -5.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.205'
$ws.Range('E24').Value = '  -2.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.359'
$ws.Range('E25').Value = '  +3.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.61'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.665'
$ws.Range('E27').Value = '  -4.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.19'
$ws.Range('E28').Value = '  -4.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '109.40'
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.351'
$ws.Range('E30').Value = '  -7.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.285'
$ws.Range('E31').Value = '  -7.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08770'
$ws.Range('E32').Value = '  -3.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04854'
$ws.Range('E33').Value = '  -3.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.159'
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7280'
$ws.Range('E35').Value = '  -10.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.858'
$ws.Range('E36').Value = '  -3.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.195'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9973'
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.405'
$ws.Range('E39').Value = '  -9.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01846'
$ws.Range('E40').Value = '  -5.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5134'
$ws.Range('E41').Value = '  -15.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9451'
$ws.Range('E42').Value = '  -11.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.93'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.225'
$ws.Range('E44').Value = '  -3.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.018'
$ws.Range('E45').Value = '  -8.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9977'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1361'
$ws.Range('E47').Value = '  -8.40%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4494'
$ws.Range('E48').Value = '  -15.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.343'
$ws.Range('E49').Value = '  -6.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.13'
$ws.Range('E50').Value = '  -3.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05919'
$ws.Range('E51').Value = '  -2.32%  '
